# Wrapping up test file audit:
# The "optimization_parameters" sheet had a stray leftover row (row 16,
# labeled "Sheet") that duplicated/obsoleted the real "simulation_timepoints"
# row below it. Remove that stray row (Excel shifts everything below it up
# by one), then leave the workbook focused on the "threshold_b" sheet with
# cell A2 selected, matching the final state the file was saved in.

$wb = $excel.ActiveWorkbook

$paramsSheet = $wb.Worksheets.Item("optimization_parameters")
$paramsSheet.Activate()
$paramsSheet.Rows.Item(16).Select()
$paramsSheet.Rows.Item(16).Delete()

$thresholdSheet = $wb.Worksheets.Item("threshold_b")
$thresholdSheet.Activate()
$thresholdSheet.Range("A2").Select()
